$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = -21.846
$ws.Range("A12").Value = -21.546
$ws.Range("A18").Value = -22.035
$ws.Range("A37").Value = -19.997
$ws.Range("A55").Value = -22.184
$ws.Range("A68").Value = -21.567
$ws.Range("A77").Value = -20.666
$ws.Range("A78").Value = -19.951
